# Auto-generated: applies the cell-value updates from the commit diff.
# All target cells are plain numeric values (no formulas) on 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Some rows gain/lose a trailing
# M/N cell entirely, which is why a couple of ClearContents() calls appear.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 289.83334
$ws.Range("J4").Value = 384.25
$ws.Range("L4").Value = 384.25
$ws.Range("N4").Value = -612.25

$ws.Range("H6").Value = 427.5263
$ws.Range("I6").Value = 320
$ws.Range("K6").Value = 960
$ws.Range("M6").Value = -848

$ws.Range("H9").Value = 10005.637
$ws.Range("I9").Value = 16927.166
$ws.Range("J9").Value = 1699.8
$ws.Range("K9").Value = 16927.166
$ws.Range("L9").Value = 1699.8
$ws.Range("M9").Value = -16758.166
$ws.Range("N9").Value = -2037.8

$ws.Range("H12").Value = 14895.429
$ws.Range("I12").Value = 20671.2
$ws.Range("J12").Value = 456
$ws.Range("K12").Value = 20671.2
$ws.Range("L12").Value = 456
$ws.Range("M12").Value = -20501.2
$ws.Range("N12").Value = -796

$ws.Range("H15").Value = 503.62296
$ws.Range("I15").Value = 503.62296
$ws.Range("K15").Value = 1510.86888
$ws.Range("M15").Value = -1341.86888

$ws.Range("H28").Value = 1257.2069
$ws.Range("I28").Value = 1049.6154
$ws.Range("K28").Value = 1049.6154
$ws.Range("M28").Value = -564.6153999999999

$ws.Range("H33").Value = 212.2
$ws.Range("I33").Value = 227.1
$ws.Range("K33").Value = 227.1
$ws.Range("M33").Value = 1.900000000000006

$ws.Range("H40").Value = 3788.3157
$ws.Range("I40").Value = 2998.4666
$ws.Range("J40").Value = 6750.25
$ws.Range("K40").Value = 2998.4666
$ws.Range("L40").Value = 6750.25
$ws.Range("M40").Value = -2823.4666
$ws.Range("N40").Value = -7100.25

$ws.Range("H70").Value = 5690.4287
$ws.Range("J70").Value = 15959.8
$ws.Range("L70").Value = 47879.39999999999
$ws.Range("N70").Value = -48419.39999999999

$ws.Range("H73").Value = 5690.4287
$ws.Range("J73").Value = 15959.8
$ws.Range("L73").Value = 47879.39999999999
$ws.Range("N73").Value = -49751.39999999999

$ws.Range("H76").Value = 8113.067
$ws.Range("I76").Value = 7549.5
$ws.Range("J76").Value = 8757.143
$ws.Range("K76").Value = 7549.5
$ws.Range("L76").Value = 8757.143
$ws.Range("M76").Value = -7234.5
$ws.Range("N76").Value = -9387.143

$ws.Range("H79").Value = 8113.067
$ws.Range("I79").Value = 7549.5
$ws.Range("J79").Value = 8757.143
$ws.Range("K79").Value = 7549.5
$ws.Range("L79").Value = 8757.143
$ws.Range("M79").Value = -6457.5
$ws.Range("N79").Value = -10941.143

$ws.Range("H113").Value = 4503
$ws.Range("I113").Value = 4200
$ws.Range("J113").Value = 4806
$ws.Range("K113").Value = 4200
$ws.Range("L113").Value = 4806
$ws.Range("M113").Value = -946
$ws.Range("N113").Value = -11314

$ws.Range("H138").Value = 2624.4883
$ws.Range("J138").Value = 2853.9167
$ws.Range("L138").Value = 8561.750100000001
$ws.Range("N138").Value = -18841.7501

$ws.Range("H140").Value = 74992.5
$ws.Range("J140").Value = 74992.5
$ws.Range("L140").Value = 74992.5
$ws.Range("N140").Value = -85352.5

$ws.Range("H141").Value = 4800
$ws.Range("I141").Value = 4800
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 14400
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -9220
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 460121.25
$ws.Range("I2").Value = 668741.75
$ws.Range("K2").Value = 668741.75
$ws.Range("M2").Value = -668628.75

$ws.Range("H102").Value = 50000348
$ws.Range("I102").Value = 50000348
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 50000348
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -49998726
$ws.Range("N102").ClearContents()

$ws.Range("H116").Value = 460121.25
$ws.Range("I116").Value = 668741.75
$ws.Range("K116").Value = 668741.75
$ws.Range("M116").Value = -666447.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 460121.25
$ws.Range("I3").Value = 668741.75
$ws.Range("K3").Value = 668741.75
$ws.Range("M3").Value = -668627.75

$ws.Range("H10").Value = 295
$ws.Range("I10").Value = 295
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 295
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -155
$ws.Range("N10").ClearContents()

$ws.Range("H100").Value = 16472
$ws.Range("J100").Value = 16472
$ws.Range("L100").Value = 16472
$ws.Range("N100").Value = -18636

$ws.Range("H107").Value = 44819.176
$ws.Range("I107").Value = 1268.3889
$ws.Range("K107").Value = 1268.3889
$ws.Range("M107").Value = 651.6111000000001

$ws.Range("H110").Value = 51494.75
$ws.Range("J110").Value = 51494.75
$ws.Range("L110").Value = 51494.75
$ws.Range("N110").Value = -59674.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 10291.4
$ws.Range("I7").Value = 25073.5
$ws.Range("K7").Value = 25073.5
$ws.Range("M7").Value = -24960.5

$ws.Range("H31").Value = 8152.093
$ws.Range("I31").Value = 6253.6787
$ws.Range("K31").Value = 6253.6787
$ws.Range("M31").Value = -5958.6787

$ws.Range("H34").Value = 8152.093
$ws.Range("I34").Value = 6253.6787
$ws.Range("K34").Value = 6253.6787
$ws.Range("M34").Value = -6051.6787

$ws.Range("H41").Value = 12338.889
$ws.Range("I41").Value = 12338.889
$ws.Range("K41").Value = 12338.889
$ws.Range("M41").Value = -11910.889

$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999.5
$ws.Range("K62").Value = 4999.5
$ws.Range("M62").Value = -4375.5

$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999.5
$ws.Range("K65").Value = 24997.5
$ws.Range("M65").Value = -21877.5

$ws.Range("H134").Value = 10418268
$ws.Range("I134").Value = 12501367
$ws.Range("K134").Value = 37504101
$ws.Range("M134").Value = -37501566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1321.2858
$ws.Range("I23").Value = 40.333332
$ws.Range("K23").Value = 120.999996
$ws.Range("M23").Value = 114.000004

$ws.Range("H113").Value = 41167.56
$ws.Range("I113").Value = 84683.414
$ws.Range("J113").Value = 999.0769
$ws.Range("K113").Value = 254050.242
$ws.Range("L113").Value = 2997.2307
$ws.Range("M113").Value = -251880.242
$ws.Range("N113").Value = -7337.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3099.3333
$ws.Range("I80").Value = 2649
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2649
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1651
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 3099.3333
$ws.Range("I83").Value = 2649
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 13245
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -8253
$ws.Range("N83").Value = -29984

$ws.Range("H102").Value = 4059.3333
$ws.Range("J102").Value = 25149.5
$ws.Range("L102").Value = 25149.5
$ws.Range("N102").Value = -28393.5

$ws.Range("H123").Value = 49999
$ws.Range("J123").Value = 49999
$ws.Range("L123").Value = 49999
$ws.Range("N123").Value = -54899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 947.86664
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H100").Value = 6240450
$ws.Range("I100").Value = 9074518
$ws.Range("J100").Value = 5499.9
$ws.Range("K100").Value = 9074518
$ws.Range("L100").Value = 5499.9
$ws.Range("M100").Value = -9073977
$ws.Range("N100").Value = -6581.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3077.9375
$ws.Range("I96").Value = 889.2
$ws.Range("J96").Value = 3483.2593
$ws.Range("K96").Value = 889.2
$ws.Range("L96").Value = 3483.2593
$ws.Range("M96").Value = 483.8
$ws.Range("N96").Value = -6229.2593

$ws.Range("H132").Value = 10873415
$ws.Range("I132").Value = 14287112
$ws.Range("K132").Value = 42861336
$ws.Range("M132").Value = -42858806

$ws.Range("H136").Value = 35715012
$ws.Range("I136").Value = 38461892
$ws.Range("K136").Value = 115385676
$ws.Range("M136").Value = -115383126

Write-Output "applied all Spriggan_Profits cell updates"